$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("en")
$ws2 = $wb.Worksheets.Item("es")

# --- Fix existing microclimate / mushroom intro rows (179-185) on the "en" sheet ---
$ws.Range("A179").Value = "intro_microclimate_01"
$ws.Range("B179").Value = "Looks like we landed in a peculiar place where the weather is a little bit colder than usual."
$ws.Range("A180").Value = "intro_microclimate_02"
$ws.Range("B180").Value = "This is known as a microclimate, where the local atmospheric condition can differ from the surrounding areas."
$ws.Range("A181").Value = "intro_microclimate_03"
$ws.Range("B181").Value = "In our case, we are below the slope of a hill that obscures most of the sunlight, and precipitations linger around longer."
$ws.Range("A182").Value = "mushroom_intro_01"
$ws.Range("B182").Value = "Uh oh, a mushroom started to grow near one of our plants!"
$ws.Range("A183").Value = "mushroom_intro_02"
$ws.Range("B183").Value = "Since there’s a lot of moisture in the region, the fungi that grow these mushrooms from beneath are able to absorb a lot of nutrients."
$ws.Range("A184").Value = "mushroom_intro_03"
$ws.Range("B184").Value = "Their spores appear to be harmful to all our frogs and plants, but one: the iron frog!"
$ws.Range("A185").Value = "mushroom_intro_04"
$ws.Range("B185").Value = "Just as the iron frogs can rid us of the moles, they, too, can rid us of these mushrooms."

# --- Add new level 3 intro dialog rows (186-199) on the "en" sheet ---
$ws.Range("A186").Value = "intro_climate_desert_01"
$ws.Range("B186").Value = "Oh boy, the desert climate...We sure took a wrong turn when we landed."
$ws.Range("A187").Value = "intro_climate_desert_02"
$ws.Range("B187").Value = "This dry and hot weather will leave us with little to no water for our plants…and with dry soil, not much nutrients are available."
$ws.Range("A188").Value = "intro_climate_desert_03"
$ws.Range("B188").Value = "However, there are oases nearby! This miracle of a microclimate will allow us to gather nutrients and water for our plants!"
$ws.Range("A189").Value = "hopper_intro_01"
$ws.Range("B189").Value = "Watch out, it's a grasshopper!"
$ws.Range("A190").Value = "hopper_intro_02"
$ws.Range("B190").Value = "These hopping hooligans can strive even in the hottest of climates, and they're no picky eaters either!"
$ws.Range("A191").Value = "hopper_intro_03"
$ws.Range("B191").Value = "These ones in particular have powerful hind legs that emit a dangerous wave that harms our frogs!"
$ws.Range("A192").Value = "hopper_intro_04"
$ws.Range("B192").Value = "Make sure to deploy any of these frogs to stop them on their tracks!"
$ws.Range("A193").Value = "antlion_intro_01"
$ws.Range("B193").Value = "An antlion larva has been spotted!"
$ws.Range("A194").Value = "antlion_intro_02"
$ws.Range("B194").Value = "These insects are commonly found in dry and sandy climates where they can easily excavate their pits."
$ws.Range("A195").Value = "antlion_intro_03"
$ws.Range("B195").Value = "Our poor unsuspecting frogs will fall prey to these pits if left alone."
$ws.Range("A196").Value = "antlion_intro_04"
$ws.Range("B196").Value = "Deploy the iron frog to plug away these troublesome pits for good!"
$ws.Range("A197").Value = "hazard_weather_intro_01"
$ws.Range("B197").Value = "Take cover, for there is a sand storm heading our way!"
$ws.Range("A198").Value = "hazard_weather_intro_02"
$ws.Range("B198").Value = "Every once in a while, in certain regions, a dangerous atmospheric condition occurs that can cause major destruction."
$ws.Range("A199").Value = "hazard_weather_intro_03"
$ws.Range("B199").Value = "With proper equipment and knowledge of weather pattern, one can properly prepare, and even avoid these severe weathers."

# --- Add matching key rows (186-199) on the "es" sheet (translations not yet filled in) ---
$ws2.Range("A186").Value = "intro_climate_desert_01"
$ws2.Range("A187").Value = "intro_climate_desert_02"
$ws2.Range("A188").Value = "intro_climate_desert_03"
$ws2.Range("A189").Value = "hopper_intro_01"
$ws2.Range("A190").Value = "hopper_intro_02"
$ws2.Range("A191").Value = "hopper_intro_03"
$ws2.Range("A192").Value = "hopper_intro_04"
$ws2.Range("A193").Value = "antlion_intro_01"
$ws2.Range("A194").Value = "antlion_intro_02"
$ws2.Range("A195").Value = "antlion_intro_03"
$ws2.Range("A196").Value = "antlion_intro_04"
$ws2.Range("A197").Value = "hazard_weather_intro_01"
$ws2.Range("A198").Value = "hazard_weather_intro_02"
$ws2.Range("A199").Value = "hazard_weather_intro_03"

# --- Restore final selection / active cell as seen in the authored workbook ---
$ws.Activate()
$ws.Range("B199").Select()
$ws2.Activate()
$ws2.Rows.Item(200).Select()
$ws.Activate()
